$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header row column width for column A
$ws.Columns.Item(1).ColumnWidth = 12.928135871887207

# Insert two new rows: new row 2 (aaaaaaaAAAA) and push existing BOAZ123 row down,
# then append a new row 4 (abc/coupon) after it.
$ws.Rows.Item(2).Insert()

# Row 2 - new coupon "aaaaaaaAAAA"
$ws.Cells.Item(2, 1).Value = "aaaaaaaAAAA"
$ws.Cells.Item(2, 2).Value = "string"
$ws.Cells.Item(2, 3).Value = "13/11/2024"
$ws.Cells.Item(2, 4).Value = "13/11/2024"
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = "Yes"
$ws.Cells.Item(2, 7).Value = "Yes"
$ws.Cells.Item(2, 8).Value = 0

# Row 3 - existing BOAZ123 coupon, CreationDateTime updated to 13/11/2024
$ws.Cells.Item(3, 1).Value = "BOAZ123"
$ws.Cells.Item(3, 2).Value = "ttt"
$ws.Cells.Item(3, 3).Value = "13/11/2024"
$ws.Cells.Item(3, 4).Value = "21/11/2024"
$ws.Cells.Item(3, 5).Value = 10
$ws.Cells.Item(3, 6).Value = "Yes"
$ws.Cells.Item(3, 7).Value = "Yes"
$ws.Cells.Item(3, 8).Value = 20

# Row 4 - new coupon "abc"
$ws.Cells.Item(4, 1).Value = "abc"
$ws.Cells.Item(4, 2).Value = "coupon"
$ws.Cells.Item(4, 3).Value = "13/11/2024"
$ws.Cells.Item(4, 4).Value = "13/11/2024"
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = "Yes"
$ws.Cells.Item(4, 7).Value = "Yes"
$ws.Cells.Item(4, 8).Value = 0
